# PSQ Floor Percentage - Create, Search, Edit and Delete test data
# Adds two new worksheets (ReleaseOrder, PSQFloorPercentage) to the
# OETestData workbook, modelled on the existing OELogin sheet.

$wb = $excel.ActiveWorkbook
$oeLogin = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create the two new worksheets, in order, right after OELogin
# ---------------------------------------------------------------------
$releaseOrder = $wb.Worksheets.Add($null, $oeLogin)
$releaseOrder.Name = "ReleaseOrder"

$psqFloor = $wb.Worksheets.Add($null, $releaseOrder)
$psqFloor.Name = "PSQFloorPercentage"

# ---------------------------------------------------------------------
# Helper seed cell: build the "bordered, general-format" style (used for
# plain data cells) once, on a scratch cell, then reuse it everywhere by
# copying formats across - this keeps styles.xml minimal/reused instead
# of minting a new xf every time.
# ---------------------------------------------------------------------
function New-BorderOnlySeed($ws, $addr) {
    $ws.Range($addr).Value = "seed"
    $ws.Range($addr).Borders.Color = 0
    $ws.Range($addr).Borders.LineStyle = 1
}

New-BorderOnlySeed $releaseOrder "Z1"
$releaseOrder.Range("Z1").Copy()
$borderOnlySeed = $releaseOrder.Range("Z1")

# ---------------------------------------------------------------------
# 2. ReleaseOrder sheet content
# ---------------------------------------------------------------------

# hrms_id / password columns - copy formatting from OELogin (keeps s=3
# header style, s=2 id style, s=1 hyperlink style, and the hyperlink
# itself all identical to the source sheet).
$releaseOrder.Range("B2").Value = "Pass@123"
$releaseOrder.Hyperlinks.Add($releaseOrder.Range("B2"), "mailto:Pass@123")

$oeLogin.Range("A1:B2").Copy()
$releaseOrder.Range("A1:B2").PasteSpecial(-4122)

$releaseOrder.Range("A1").Value = "hrms_id"
$releaseOrder.Range("B1").Value = "password"
$releaseOrder.Range("A2").Value = "110"

# Header row (C1:H1) - same bold/filled header style as A1/B1
$oeLogin.Range("A1").Copy()
$releaseOrder.Range("C1:H1").PasteSpecial(-4122)
$releaseOrder.Range("C1").Value = "StoreCode"
$releaseOrder.Range("D1").Value = "OrderNo"
$releaseOrder.Range("E1").Value = "OrderType"
$releaseOrder.Range("F1").Value = "Vendor"
$releaseOrder.Range("G1").Value = "Segment"
$releaseOrder.Range("H1").Value = "OrderKind"

# Data row (C2:H2): C2/F2/G2/H2 -> border-only style; D2/E2 -> text style
$borderOnlySeed.Copy()
$releaseOrder.Range("C2").PasteSpecial(-4122)
$releaseOrder.Range("C2").Value = "TS Madision"

$oeLogin.Range("A2").Copy()
$releaseOrder.Range("D2:E2").PasteSpecial(-4122)
$releaseOrder.Range("D2").Value = "48608"
$releaseOrder.Range("E2").Value = "Customer Order"

$borderOnlySeed.Copy()
$releaseOrder.Range("F2:H2").PasteSpecial(-4122)
$releaseOrder.Range("F2").Value = "BKK - BALKISHORE KHANNA AND COMPANY"
$releaseOrder.Range("G2").Value = "Gold"
$releaseOrder.Range("H2").Value = "New Order-NO"

$releaseOrder.Range("Z1").Clear()
$releaseOrder.Range("A1:B2").Select()

# ---------------------------------------------------------------------
# 3. PSQFloorPercentage sheet content
# ---------------------------------------------------------------------
New-BorderOnlySeed $psqFloor "Z1"
$psqFloor.Range("Z1").Copy()
$borderOnlySeed2 = $psqFloor.Range("Z1")

$psqFloor.Range("B2").Value = "Pass@123"
$psqFloor.Hyperlinks.Add($psqFloor.Range("B2"), "mailto:Pass@123")

$oeLogin.Range("A1:B2").Copy()
$psqFloor.Range("A1:B2").PasteSpecial(-4122)

$psqFloor.Range("A1").Value = "hrms_id"
$psqFloor.Range("B1").Value = "password"
$psqFloor.Range("A2").Value = "110"

# Header row (C1:M1)
$oeLogin.Range("A1").Copy()
$psqFloor.Range("C1:M1").PasteSpecial(-4122)
$psqFloor.Range("C1").Value = "company"
$psqFloor.Range("D1").Value = "region"
$psqFloor.Range("E1").Value = "role"
$psqFloor.Range("F1").Value = "grade"
$psqFloor.Range("G1").Value = "storename"
$psqFloor.Range("H1").Value = "segment"
$psqFloor.Range("I1").Value = "floorpercentage"
$psqFloor.Range("J1").Value = "year"
$psqFloor.Range("K1").Value = "month"
$psqFloor.Range("L1").Value = "is_active"
$psqFloor.Range("M1").Value = "editfloorpercentagevalue"

# Data row (C2:M2)
$oeLogin.Range("A2").Copy()
$psqFloor.Range("C2").PasteSpecial(-4122)
$psqFloor.Range("C2").Value = "CKCCO-C. Krishniah Chetty & Co. Private Limited "

$borderOnlySeed2.Copy()
$psqFloor.Range("D2:H2").PasteSpecial(-4122)
$psqFloor.Range("D2").Value = "Bangalore"
$psqFloor.Range("E2").Value = "Valuator"
$psqFloor.Range("F2").Value = "A"
$psqFloor.Range("G2").Value = "TS Madision"
$psqFloor.Range("H2").Value = "Gold"

$oeLogin.Range("A2").Copy()
$psqFloor.Range("I2:J2").PasteSpecial(-4122)
$psqFloor.Range("I2").Value = "50"
$psqFloor.Range("J2").Value = "2023"

$borderOnlySeed2.Copy()
$psqFloor.Range("K2").PasteSpecial(-4122)
$psqFloor.Range("K2").Value = "Jul"

$oeLogin.Range("A2").Copy()
$psqFloor.Range("L2:M2").PasteSpecial(-4122)
$psqFloor.Range("L2").Value = "1"
$psqFloor.Range("M2").Value = "75"

$psqFloor.Range("Z1").Clear()

# ---------------------------------------------------------------------
# 4. Column widths (approximate auto-fit sizing used by the author)
# ---------------------------------------------------------------------
$releaseOrder.Columns.Item(1).ColumnWidth = 10
$releaseOrder.Columns.Item(2).ColumnWidth = 7.5
$releaseOrder.Columns.Item(3).ColumnWidth = 14.5
$releaseOrder.Columns.Item(4).ColumnWidth = 36.5
$releaseOrder.Columns.Item(5).ColumnWidth = 13.5
$releaseOrder.Columns.Item(6).ColumnWidth = 12.5
$releaseOrder.Columns.Item(8).ColumnWidth = 12.5

$psqFloor.Columns.Item(3).ColumnWidth = 40.5
$psqFloor.Columns.Item(4).ColumnWidth = 8
$psqFloor.Columns.Item(5).ColumnWidth = 8
$psqFloor.Columns.Item(7).ColumnWidth = 11.5
$psqFloor.Columns.Item(9).ColumnWidth = 13.5
$psqFloor.Columns.Item(13).ColumnWidth = 21.5

# ---------------------------------------------------------------------
# 5. Sheet view / selection state - OELogin loses its "tabSelected" flag
#    (it moves to the background) and its leftover M19 selection is reset.
# ---------------------------------------------------------------------
$oeLogin.Range("A1:B2").Select()

# Make PSQFloorPercentage (the last created sheet) the active / visible
# tab, matching activeTab=2 in the saved workbook view, with its
# selection parked on L8.
$psqFloor.Select()
$psqFloor.Range("L8").Select()
